$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Aula 10 - "Importando os arquivos HTML" - new rows 4 and 5.
# Set the text cells first (in the order the shared-string table should
# grow) before touching row 3, so the rebuilt sharedStrings.xml lines up.
$ws.Range("C4").Value = "2. Inciando o Desenvolvimento do Projeto"
$ws.Range("D4").Value = "10. Importanto os arquivos HTML"

# Row 3: fix the "sessão" text for aula 9.
$ws.Range("C3").Value = "2. Iniciando o desenvolvimento do projeto"

# Row 5's observação note, then row 4's observação note.
$ws.Range("E5").Value = "2:53`npor padrão, as paginas HTML devem ficar no diretório TEMPLATES pois é lá que o spring MVC procura as páginas. Caso deseja alterar esse diretório padrão, é necessário sobrescrever o arquivo de propriedades do thymeleaf (professor citou que isso será visto em aulas posteriores)"
$ws.Range("E4").Value = "4:35`nO spring MVC só encontra as páginas HTML através de um CONTROLLER`n"

# Remaining cells for the two new rows.
$ws.Range("B4").Value = 10
$ws.Range("C5").Value = "2. Inciando o Desenvolvimento do Projeto"
$ws.Range("D5").Value = "10. Importanto os arquivos HTML"
$ws.Range("B5").Value = 10

# Wrap text + row heights for the two new note cells.
$ws.Range("E4").WrapText = $true
$ws.Range("E5").WrapText = $true
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 75

$ws.Range("D9").Select()
